$d = $word.ActiveDocument

# --- 1) Update the two AbstractTemplatesTestSuite line numbers and the
#        GeneratedMethodAccessor index, a few lines above the big block
#        that gets rewritten below.
$null = $d.Content.Find.Execute(
    "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)",
    2)

$null = $d.Content.Find.Execute(
    "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)",
    2)

$null = $d.Content.Find.Execute(
    "GeneratedMethodAccessor5",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "GeneratedMethodAccessor4",
    2)

# --- 2) Replace the tail of the stack trace: everything from the
#        "JUnit4Provider.execute(...)" line through the final
#        "Main.main(...)" line gets replaced by the new
#        Eclipse JDT-runner stack frames.
$startRange = $d.Content.Duplicate
$null = $startRange.Find.Execute(
    "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$endRange = $d.Content.Duplicate
$null = $endRange.Find.Execute(
    "at org.eclipse.equinox.launcher.Main.main(Main.java:1420)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target = $d.Range($startRange.Start, $endRange.End)
$newTail = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n" +
           "`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n" +
           "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n" +
           "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n" +
           "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n" +
           "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
$target.Text = $newTail

Write-Output "edit complete"
